# Update UUT (E column) measured values across the voltage/current/phase
# calibration tables to reflect the new "high voltage and current" test run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24.9979763
$ws.Range("E3").Value = 49.99749374
$ws.Range("E4").Value = 74.99463654
$ws.Range("E5").Value = 99.99442291
$ws.Range("E6").Value = 124.99074554
$ws.Range("E7").Value = 149.98817444
$ws.Range("E8").Value = 25.01117325
$ws.Range("E9").Value = 50.02285767
$ws.Range("E10").Value = 75.03379821999999
$ws.Range("E11").Value = 100.04698944
$ws.Range("E12").Value = 125.06109619
$ws.Range("E13").Value = 150.07202148
$ws.Range("E14").Value = 24.99144936
$ws.Range("E15").Value = 49.98563385
$ws.Range("E16").Value = 74.97724915000001
$ws.Range("E17").Value = 99.97194672000001
$ws.Range("E18").Value = 124.96549988
$ws.Range("E19").Value = 149.95675659
$ws.Range("E26").Value = 25.00974274
$ws.Range("E27").Value = 50.01837158
$ws.Range("E28").Value = 75.02858734
$ws.Range("E29").Value = 100.03844452
$ws.Range("E30").Value = 125.05028534
$ws.Range("E31").Value = 150.06365967
$ws.Range("E32").Value = 25.00335503
$ws.Range("E33").Value = 50.00356293
$ws.Range("E34").Value = 75.00463867000001
$ws.Range("E35").Value = 100.01135254
$ws.Range("E36").Value = 125.01959991
$ws.Range("E37").Value = 150.02742004
$ws.Range("E74").Value = 60.26304637721221
$ws.Range("E75").Value = 120.2972119614733
$ws.Range("E76").Value = 179.7593575491723
$ws.Range("E77").Value = 60.23295852263566
$ws.Range("E78").Value = 120.2331134184667
$ws.Range("E79").Value = 179.768394454297
$ws.Range("E80").Value = 60.21943581112865
$ws.Range("E81").Value = 120.2084549933066
$ws.Range("E82").Value = 179.8082515568856
$ws.Range("E83").Value = 60.24495290608287
$ws.Range("E84").Value = 120.2447823746622
$ws.Range("E86").Value = 60.20701054828241
$ws.Range("E87").Value = 120.167535693746
$ws.Range("E89").Value = 60.20570264038873
$ws.Range("E90").Value = 119.8126910889515
